$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "FUL"/"Q" record and the "RDB"/"Z" record were being dropped from the
# dataset when the grid was re-saved; re-enter the full, correctly ordered
# dataset (rows 8-21) so both records are preserved and land back in the
# table - "FUL"/"Q" after "LUB"/"P" and "RDB"/"Z" at the very end - while
# also correcting a typo ("RFD" -> "RDF").
$ws.Cells.Item(8, 1).Value = "FDL"
$ws.Cells.Item(8, 2).Value = "I"

$ws.Cells.Item(9, 1).Value = "FDR"
$ws.Cells.Item(9, 2).Value = "K"

$ws.Cells.Item(10, 1).Value = "BUR"
$ws.Cells.Item(10, 2).Value = "L"

$ws.Cells.Item(11, 1).Value = "BUL"
$ws.Cells.Item(11, 2).Value = "M"

$ws.Cells.Item(12, 1).Value = "BDR"
$ws.Cells.Item(12, 2).Value = "N"

$ws.Cells.Item(13, 1).Value = "BDL"
$ws.Cells.Item(13, 2).Value = "O"

$ws.Cells.Item(14, 1).Value = "LUB"
$ws.Cells.Item(14, 2).Value = "P"

$ws.Cells.Item(15, 1).Value = "FUL"
$ws.Cells.Item(15, 2).Value = "Q"

$ws.Cells.Item(16, 1).Value = "LDB"
$ws.Cells.Item(16, 2).Value = "R"

$ws.Cells.Item(17, 1).Value = "LUF"
$ws.Cells.Item(17, 2).Value = "S"

$ws.Cells.Item(18, 1).Value = "LDF"
$ws.Cells.Item(18, 2).Value = "T"

$ws.Cells.Item(19, 1).Value = "RDF"
$ws.Cells.Item(19, 2).Value = "U"

$ws.Cells.Item(20, 1).Value = "RUB"
$ws.Cells.Item(20, 2).Value = "W"

$ws.Cells.Item(21, 1).Value = "RDB"
$ws.Cells.Item(21, 2).Value = "Z"

# Excel re-wrapped a handful of rows to a slightly shorter height after the
# dataset was rebuilt.
$ws.Rows("7").RowHeight = 18.75
$ws.Rows("8").RowHeight = 18.75
$ws.Rows("11").RowHeight = 18.75
$ws.Rows("15").RowHeight = 18.75
